# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF column (F) values to reflect the repulled data
$ws.Range("F2").Value = -2
$ws.Range("F5").Value = -2
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = -5
$ws.Range("F8").Value = -7
$ws.Range("F10").Value = -3
